# Refresh cryptocurrency price (Price) and 1h volume change (Volume(1h)) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "60.916.90"
$ws.Cells.Item(2, 5).Value2 = "  -0.74%  "
$ws.Cells.Item(3, 4).Value2 = "3.361.88"
$ws.Cells.Item(3, 5).Value2 = "  -1.23%  "
$ws.Cells.Item(4, 5).Value2 = "  +0.04%  "
$ws.Cells.Item(5, 4).Value2 = "405.87"
$ws.Cells.Item(5, 5).Value2 = "  -1.58%  "
$ws.Cells.Item(6, 4).Value2 = "135.26"
$ws.Cells.Item(6, 5).Value2 = "  +10.69%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = "0.590"
$ws.Cells.Item(7, 5).Value2 = "  +2.49%  "
$ws.Cells.Item(8, 5).Value2 = "  +0.05%  "
$ws.Cells.Item(9, 4).Value2 = "0.674"
$ws.Cells.Item(9, 5).Value2 = "  +5.78%  "
$ws.Cells.Item(10, 5).Value2 = "  +3.74%  "
$ws.Cells.Item(11, 4).Value2 = "42.54"
$ws.Cells.Item(11, 5).Value2 = "  +3.70%  "
$ws.Cells.Item(12, 5).Value2 = "  -0.86%  "
$ws.Cells.Item(13, 4).Value2 = "3.892.14"
$ws.Cells.Item(13, 5).Value2 = "  -1.46%  "
$ws.Cells.Item(14, 4).Value2 = "8.32"
$ws.Cells.Item(14, 5).Value2 = "  -0.71%  "
$ws.Cells.Item(15, 4).Value2 = "19.58"
$ws.Cells.Item(15, 5).Value2 = "  +0.56%  "
$ws.Cells.Item(16, 4).Value2 = "3.357.74"
$ws.Cells.Item(16, 5).Value2 = "  -1.02%  "
$ws.Cells.Item(17, 4).Value2 = "60.945.01"
$ws.Cells.Item(17, 5).Value2 = "  -0.70%  "
$ws.Cells.Item(18, 4).Value2 = "1.01"
$ws.Cells.Item(18, 5).Value2 = "  -0.25%  "
$ws.Cells.Item(19, 4).Value2 = "10.97"
$ws.Cells.Item(19, 5).Value2 = "  +1.55%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = "0.0000127"
$ws.Cells.Item(20, 5).Value2 = "  +5.11%  "
$ws.Cells.Item(21, 4).Value2 = "3.21"
$ws.Cells.Item(21, 5).Value2 = "  -3.28%  "
$ws.Cells.Item(22, 4).Value2 = "83.76"
$ws.Cells.Item(22, 5).Value2 = "  +10.08%  "
$ws.Cells.Item(23, 4).Value2 = "308.66"
$ws.Cells.Item(23, 5).Value2 = "  +3.58%  "
$ws.Cells.Item(24, 4).Value2 = "12.66"
$ws.Cells.Item(24, 5).Value2 = "  -1.05%  "
$ws.Cells.Item(25, 4).Value2 = "3.12"
$ws.Cells.Item(25, 5).Value2 = "  -0.38%  "
$ws.Cells.Item(26, 5).Value2 = "  +12.08%  "
$ws.Cells.Item(27, 4).Value2 = "8.34"
$ws.Cells.Item(27, 5).Value2 = "  +9.51%  "
$ws.Cells.Item(28, 4).Value2 = "29.37"
$ws.Cells.Item(28, 5).Value2 = "  -4.08%  "
$ws.Cells.Item(29, 4).Value2 = "7.43"
$ws.Cells.Item(29, 5).Value2 = "  -8.22%  "
$ws.Cells.Item(30, 4).Value2 = "0.172"
$ws.Cells.Item(30, 5).Value2 = "  +0.26%  "
$ws.Cells.Item(31, 4).Value2 = "0.116"
$ws.Cells.Item(31, 5).Value2 = "  +0.68%  "
$ws.Cells.Item(32, 5).Value2 = "  -0.08%  "
$ws.Cells.Item(33, 4).Value2 = "11.27"
$ws.Cells.Item(33, 5).Value2 = "  -0.80%  "
$ws.Cells.Item(34, 4).Value2 = "41.16"
$ws.Cells.Item(34, 5).Value2 = "  -2.97%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = "2.50"
$ws.Cells.Item(35, 5).Value2 = "  -1.92%  "
$ws.Cells.Item(36, 4).Value2 = "0.0479"
$ws.Cells.Item(36, 5).Value2 = "  +0.04%  "
$ws.Cells.Item(37, 4).Value2 = "51.83"
$ws.Cells.Item(37, 5).Value2 = "  -1.49%  "
$ws.Cells.Item(38, 5).Value2 = "  -0.10%  "
$ws.Cells.Item(39, 4).Value2 = "3.41"
$ws.Cells.Item(39, 5).Value2 = "  -3.34%  "
$ws.Cells.Item(40, 5).Value2 = "  -3.32%  "
$ws.Cells.Item(41, 4).Value2 = "1.98"
$ws.Cells.Item(41, 5).Value2 = "  +0.70%  "
$ws.Cells.Item(42, 4).Value2 = "136.96"
$ws.Cells.Item(42, 5).Value2 = "  +2.74%  "
$ws.Cells.Item(43, 4).Value2 = "0.123"
$ws.Cells.Item(43, 5).Value2 = "  +0.86%  "
$ws.Cells.Item(44, 5).Value2 = "  +2.58%  "
$ws.Cells.Item(45, 4).Value2 = "0.286"
$ws.Cells.Item(45, 5).Value2 = "  +1.51%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = "16.60"
$ws.Cells.Item(46, 5).Value2 = "  -3.18%  "
$ws.Cells.Item(47, 5).Value2 = "  +1.39%  "
$ws.Cells.Item(48, 4).Value2 = "21.42"
$ws.Cells.Item(48, 5).Value2 = "  -1.24%  "
$ws.Cells.Item(49, 4).Value2 = "2.119.02"
$ws.Cells.Item(49, 5).Value2 = "  -3.79%  "
$ws.Cells.Item(50, 5).Value2 = "  -4.14%  "
$ws.Cells.Item(51, 5).Value2 = "  -2.17%  "
